$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow the default/body columns (A, C:H) and shrink the three
# custom-width columns (B, I, J) to build the new, more compact CDE view.
$ws.Columns.Item(1).ColumnWidth = 7.671768707482998
$ws.Range("C1:H1").EntireColumn.ColumnWidth = 7.671768707482998
$ws.Columns.Item(2).ColumnWidth = 39.12585034013606
$ws.Columns.Item(9).ColumnWidth = 26.97278911564627
$ws.Columns.Item(10).ColumnWidth = 15.636054421768668

# Move the selection onto the newly relevant CDE cells.
[void]$ws.Range("I24:I25").Select()
